$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two percentage columns (F = "Errors requiring Analysis %",
# D = "Timeout Errors %"). Delete right-to-left so column letters stay
# valid for the second delete.
$ws.Columns("F").Delete()
$ws.Columns("D").Delete()

# Update the last data row with the latest counts (row 10), now that the
# layout is A:Date B:Total C:SessionTimeout D:ErrorsRequiringAnalysis.
$ws.Range("B10").Value = 652
$ws.Range("C10").Value = 41
$ws.Range("D10").Value = 611

# Tidy column A width to match the new layout.
$ws.Columns("A").ColumnWidth = 10.14

# Restore the cursor/selection shown when the file was last saved.
$ws.Range("G10").Select() | Out-Null
